$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Entities")

# Update the magicAtk value for row 4 (10103 / Boss) from 0.7 to 0.8
$ws.Range("I4").Value = 0.8

# Move the active selection to I5 (mirrors the saved selection state in the file)
$ws.Range("I5").Select()
